# Insert a new row at position 375 (shifting existing rows 375-439 down to 376-440)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(375).Insert()

# Populate the newly inserted row 375 with its data
$ws.Range("A375").Value = 10
$ws.Range("B375").Value = "Vega Modelo de Temuco"
$ws.Range("C375").Value = "La Araucanía"
$ws.Range("D375").Value = 44694
$ws.Range("E375").Value = 9
$ws.Range("F375").Value = 100112032
$ws.Range("G375").Value = "Zapallo italiano"
$ws.Range("H375").Value = "Sin especificar"
$ws.Range("I375").Value = "Primera"
$ws.Range("J375").Value = 125
$ws.Range("K375").Value = 20000
$ws.Range("L375").Value = 20000
$ws.Range("M375").Value = 20000
$ws.Range("N375").Value = "$/caja 60 unidades"
$ws.Range("O375").Value = "Región de Arica y Parinacota"
$ws.Range("P375").Value = 333
$ws.Range("Q375").Value = 60
$ws.Range("R375").Value = "Hortaliza"
